# Swap the presentation's theme palette with the Notes Master's theme
# palette: the main deck (Slide Master / Presentation theme, stored as
# ppt/theme/theme1.xml) moves from the custom "Integral" palette to the
# default "Office Theme" palette (the one the Notes Master previously
# used, stored as ppt/theme/theme2.xml).
#
# fontScheme / fmtScheme are identical between the two themes already,
# so only the 12 theme colours (clrScheme) need to change.

$p = $ppt.ActivePresentation

$master = $p.SlideMaster
$clrScheme = $master.Theme.ThemeColorScheme

# Office Theme colour values (RRGGBB), applied in the standard
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order (indices 1-12).
$clrScheme.Item(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1
$clrScheme.Item(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1
$clrScheme.Item(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2
$clrScheme.Item(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2
$clrScheme.Item(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1
$clrScheme.Item(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2
$clrScheme.Item(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3
$clrScheme.Item(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4
$clrScheme.Item(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5
$clrScheme.Item(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6
$clrScheme.Item(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink
$clrScheme.Item(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink
